# Auto-generated edit script: updates computed profit/price columns (H-N)
# across multiple worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 528.1
$ws.Range("I2").Value = 374.75
$ws.Range("J2").Value = 630.3333
$ws.Range("K2").Value = 374.75
$ws.Range("L2").Value = 630.3333
$ws.Range("M2").Value = -261.75
$ws.Range("N2").Value = -856.3333
# Row 9
$ws.Range("H9").Value = 433.33334
$ws.Range("I9").Value = 150
$ws.Range("J9").Value = 575
$ws.Range("K9").Value = 150
$ws.Range("L9").Value = 575
$ws.Range("M9").Value = 19
$ws.Range("N9").Value = -913
# Row 15
$ws.Range("H15").Value = 205.64
$ws.Range("I15").Value = 205.64
$ws.Range("K15").Value = 616.92
$ws.Range("M15").Value = -447.92
# Row 17
$ws.Range("H17").Value = 2487
$ws.Range("J17").Value = 2516.4
$ws.Range("L17").Value = 7549.200000000001
$ws.Range("N17").Value = -7885.200000000001
# Row 29
$ws.Range("H29").Value = 202.5
$ws.Range("I29").Value = 202.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 607.5
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -326.5
$ws.Range("N29").ClearContents()
# Row 76
$ws.Range("H76").Value = 3971196
$ws.Range("I76").Value = 3080
$ws.Range("J76").Value = 6175705
$ws.Range("K76").Value = 3080
$ws.Range("L76").Value = 6175705
$ws.Range("M76").Value = -2765
$ws.Range("N76").Value = -6176335
# Row 79
$ws.Range("H79").Value = 3971196
$ws.Range("I79").Value = 3080
$ws.Range("J79").Value = 6175705
$ws.Range("K79").Value = 3080
$ws.Range("L79").Value = 6175705
$ws.Range("M79").Value = -1988
$ws.Range("N79").Value = -6177889
# Row 94
$ws.Range("H94").Value = 3000
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -2549
$ws.Range("N94").ClearContents()
# Row 116
$ws.Range("H116").Value = 19236408
$ws.Range("I116").Value = 83335370
$ws.Range("J116").Value = 6720.3
$ws.Range("K116").Value = 83335370
$ws.Range("L116").Value = 6720.3
$ws.Range("M116").Value = -83331928
$ws.Range("N116").Value = -13604.3
# Row 129
$ws.Range("H129").Value = 752.89795
$ws.Range("J129").Value = 800.11365
$ws.Range("L129").Value = 2400.34095
$ws.Range("N129").Value = -12400.34095
# Row 132
$ws.Range("H132").Value = 3256.25
$ws.Range("I132").Value = 3228.6072
$ws.Range("K132").Value = 9685.821599999999
$ws.Range("M132").Value = -7155.821599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6969.046
$ws.Range("I32").Value = 5756.4478
$ws.Range("J32").Value = 11031.25
$ws.Range("K32").Value = 5756.4478
$ws.Range("L32").Value = 11031.25
$ws.Range("M32").Value = -5469.4478
$ws.Range("N32").Value = -11605.25
# Row 74
$ws.Range("H74").Value = 28572880
$ws.Range("I74").Value = 35714908
$ws.Range("J74").Value = 4771.4287
$ws.Range("K74").Value = 35714908
$ws.Range("L74").Value = 4771.4287
$ws.Range("M74").Value = -35714034
$ws.Range("N74").Value = -6519.4287
# Row 77
$ws.Range("H77").Value = 28572880
$ws.Range("I77").Value = 35714908
$ws.Range("J77").Value = 4771.4287
$ws.Range("K77").Value = 178574540
$ws.Range("L77").Value = 23857.1435
$ws.Range("M77").Value = -178570172
$ws.Range("N77").Value = -32593.1435
# Row 97
$ws.Range("H97").Value = 1118.2
$ws.Range("I97").Value = 1118.2
$ws.Range("K97").Value = 1118.2
$ws.Range("M97").Value = -622.2

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4369.227
$ws.Range("I31").Value = 2370.7646
$ws.Range("J31").Value = 5627.5186
$ws.Range("K31").Value = 2370.7646
$ws.Range("L31").Value = 5627.5186
$ws.Range("M31").Value = -2075.7646
$ws.Range("N31").Value = -6217.5186
# Row 34
$ws.Range("H34").Value = 4369.227
$ws.Range("I34").Value = 2370.7646
$ws.Range("J34").Value = 5627.5186
$ws.Range("K34").Value = 2370.7646
$ws.Range("L34").Value = 5627.5186
$ws.Range("M34").Value = -2168.7646
$ws.Range("N34").Value = -6031.5186
# Row 94
$ws.Range("H94").Value = 3433
$ws.Range("I94").Value = 2772
$ws.Range("J94").Value = 3895.7
$ws.Range("K94").Value = 2772
$ws.Range("L94").Value = 3895.7
$ws.Range("M94").Value = -2321
$ws.Range("N94").Value = -4797.7
# Row 122
$ws.Range("H122").Value = 1328.8422
$ws.Range("I122").Value = 956.1667
$ws.Range("K122").Value = 2868.5001
$ws.Range("M122").Value = -418.5001000000002
# Row 132
$ws.Range("H132").Value = 4035.4443
$ws.Range("I132").Value = 3113.6
$ws.Range("J132").Value = 5187.75
$ws.Range("K132").Value = 9340.799999999999
$ws.Range("L132").Value = 15563.25
$ws.Range("M132").Value = -6810.799999999999
$ws.Range("N132").Value = -20623.25

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 755.1900000000001
$ws.Range("J131").Value = 767.3299
$ws.Range("L131").Value = 2301.9897
$ws.Range("N131").Value = -12381.9897
# Row 140
$ws.Range("H140").Value = 3154.1428
$ws.Range("I140").Value = 1973.8889
$ws.Range("J140").Value = 5278.6
$ws.Range("K140").Value = 5921.6667
$ws.Range("L140").Value = 15835.8
$ws.Range("M140").Value = -741.6666999999998
$ws.Range("N140").Value = -26195.8

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1358.4482
$ws.Range("I97").Value = 1583.75
$ws.Range("J97").Value = 857.7778
$ws.Range("K97").Value = 1583.75
$ws.Range("L97").Value = 857.7778
$ws.Range("M97").Value = -1087.75
$ws.Range("N97").Value = -1849.7778
# Row 102
$ws.Range("H102").Value = 2854.1
$ws.Range("I102").Value = 2414.889
$ws.Range("J102").Value = 6807
$ws.Range("K102").Value = 2414.889
$ws.Range("L102").Value = 6807
$ws.Range("M102").Value = -792.8890000000001
$ws.Range("N102").Value = -10051

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 1229214
$ws.Range("I122").Value = 1786247.6
$ws.Range("K122").Value = 5358742.800000001
$ws.Range("M122").Value = -5356292.800000001
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1474.375
$ws.Range("I122").Value = 1050
$ws.Range("K122").Value = 3150
$ws.Range("M122").Value = -700
# Row 132
$ws.Range("H132").Value = 2034.3684
$ws.Range("I132").Value = 1950
$ws.Range("K132").Value = 5850
$ws.Range("M132").Value = -3320
# Row 136
$ws.Range("H136").Value = 27901834
$ws.Range("I136").Value = 35596410
$ws.Range("K136").Value = 106789230
$ws.Range("M136").Value = -106786680
